# Generate Report for Handoff
# Update the "Latest Handoff Date/Datetime" values for the
# 2f66370b-b656-486f-b385-4db18e61c7c7.md row (row 6) across the
# Overview, zh-cn and de-de sheets to reflect the freshly generated
# handoff report timestamps.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("D6").Value = "2016-03-24 16:48:25"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E6").Value = "2016-03-24 16:48:16"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E6").Value = "2016-03-24 16:48:25"
